$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (column C) date values for rows 2-10 from 45170 (2023-09-01)
# to 45174 (2023-09-05), matching the new serial date value while keeping
# the existing date formatting.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 3).Value = 45174
}
